$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D14").Value = 0.8087
$ws.Range("E14").Value = 0.8541
$ws.Range("F14").Value = 0.9131
$ws.Range("G14").Value = 0.976

$ws.Range("D15").Value = 0.4201
$ws.Range("E15").Value = 0.4619
$ws.Range("F15").Value = 0.4914
$ws.Range("G15").Value = 0.5138

$ws.Range("D16").Value = 0.3914
$ws.Range("E16").Value = 0.4278
$ws.Range("F16").Value = 0.4576
$ws.Range("G16").Value = 0.4769

$ws.Range("D17").Value = 0.2938
$ws.Range("E17").Value = 0.3362
$ws.Range("F17").Value = 0.3569
$ws.Range("G17").Value = 0.3716

$ws.Range("B18").Value = 0.3286
$ws.Range("D18").Value = 0.3247
$ws.Range("E18").Value = 0.3659
$ws.Range("F18").Value = 0.3785
$ws.Range("G18").Value = 0.4525

$ws.Range("B29").Value = 0.4867
$ws.Range("D29").Value = 0.5199
$ws.Range("E29").Value = 0.5653
$ws.Range("F29").Value = 0.6038
$ws.Range("G29").Value = 0.6411

$ws.Range("D30").Value = 0.3787
$ws.Range("E30").Value = 0.4152
$ws.Range("F30").Value = 0.4293
$ws.Range("G30").Value = 0.5111
